$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-7, columns B-E and G.
# F (Win) column is unchanged.
$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 0.4998867070740569;  G = 6.048734245549538 }
    3 = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 0.1529057820181812;  E = 0.4998867070740569;  G = 1.145820798638228 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 3.082599426703578;   E = 246.9852506941017;   G = 254.9039648082657 }
    5 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 0.4998867070740569;  G = 6.048734245549538 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 16.98373111632243;   E = 0.4998867070740569;  G = 22.31973251085698 }
    7 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 0.4998867070740569;  G = 6.048734245549538 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}

$wb.Save()
